$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect so the model-holdings data can be
# refreshed, then restore protection afterwards.
$ws.Unprotect()

# Bump the "as of" date in the confidentiality / disclosure banner text.
$ws.Range("A37").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-29 for illustrative purposes only and are subject to change."

# Refresh the Weight (D) and Percent Change (E) columns for each holding
# with the latest model values.
$ws.Range("D2").Value = 0.03837465681813482
$ws.Range("E2").Value = -0.004667444574095736
$ws.Range("D3").Value = 0.0216926923939289
$ws.Range("E3").Value = -0.006245120999219367
$ws.Range("D4").Value = 0.01996497120184253
$ws.Range("E4").Value = -0.005165289256198302
$ws.Range("D5").Value = 0.04067027028317921
$ws.Range("E5").Value = -0.006245662734212232
$ws.Range("D6").Value = 0.03723078511283134
$ws.Range("E6").Value = -0.0007843137254902599
$ws.Range("D7").Value = 0.02100490391233432
$ws.Range("E7").Value = 0.0003875968992248513
$ws.Range("D8").Value = 0.03765576805245299
$ws.Range("E8").Value = 0.004468275245755216
$ws.Range("D9").Value = 0.02134879815313161
$ws.Range("E9").Value = -0.005216913783635313
$ws.Range("D10").Value = 0.02600645919791987
$ws.Range("E10").Value = 0.01826150474799126
$ws.Range("D11").Value = 0.02430218534043332
$ws.Range("E11").Value = -0.00335008375209378
$ws.Range("D12").Value = 0.0585611836214506
$ws.Range("E12").Value = -0.003464488987874215
$ws.Range("D13").Value = 0.02649548809779102
$ws.Range("E13").Value = -0.006865781710914565
$ws.Range("D14").Value = 0.02744271699463354
$ws.Range("E14").Value = 0.01070385987674349
$ws.Range("D15").Value = 0.03536726477452104
$ws.Range("E15").Value = 0.01117222413200425
$ws.Range("D16").Value = 0.01908488201078697
$ws.Range("E16").Value = 0.01020689655172413
$ws.Range("D17").Value = 0.03067964324663302
$ws.Range("E17").Value = -0.01556833259619639
$ws.Range("D18").Value = 0.02397164305406092
$ws.Range("E18").Value = -0.004387990762124483
$ws.Range("D19").Value = 0.1329694003056512
$ws.Range("E19").Value = -0.002680965147452974
$ws.Range("D20").Value = 0.009661007446211298
$ws.Range("E20").Value = 0.01175862514536763
$ws.Range("D21").Value = 0.01586830934784984
$ws.Range("E21").Value = 0.01019968395345483
$ws.Range("D22").Value = 0.01745002746106236
$ws.Range("E22").Value = -0.009359137055837574
$ws.Range("D23").Value = 0.01669526210239331
$ws.Range("E23").Value = -0.0007087172218286364
$ws.Range("D24").Value = 0.02164438654255428
$ws.Range("E24").Value = -0.003773991805046339
$ws.Range("D25").Value = 0.01207320626940477
$ws.Range("E25").Value = -0.01174698795180729
$ws.Range("D26").Value = 0.04366577583078059
$ws.Range("E26").Value = 0.007998657567960699
$ws.Range("D27").Value = 0.02545642380709921
$ws.Range("E27").Value = -0.0001961553550412498
$ws.Range("D28").Value = 0.04783875007469767
$ws.Range("E28").Value = -0.001691638472692159
$ws.Range("D29").Value = 0.05770128519074488
$ws.Range("E29").Value = -0.004609144542772836
$ws.Range("D30").Value = 0.01349676342654355
$ws.Range("E30").Value = -0.02202937249666226
$ws.Range("D31").Value = 0.01448497088269825
$ws.Range("E31").Value = -0.01030444964871191
$ws.Range("D32").Value = 0.04431057753227551
$ws.Range("E32").Value = -0.0005181347150259308
$ws.Range("D33").Value = 0.01682954151396726
$ws.Range("E33").Value = -0.003792667509481706
$ws.Range("E34").Value = -0.001342012537919435

$ws.Protect()
